$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated transition-probability matrix cells after simulating more games
# (values recomputed from a larger number of simulated games per starting state)
$ws.Range("B2").Value = 0.2212189616252822
$ws.Range("C2").Value = 0.5079006772009029
$ws.Range("J2").Value = 0.01580135440180587
$ws.Range("P2").Value = 0.1896162528216704
$ws.Range("S2").Value = 0.0654627539503386
$ws.Range("B3").Value = 0.008298755186721992
$ws.Range("C3").Value = 0.03319502074688797
$ws.Range("J3").Value = 0.02074688796680498
$ws.Range("P3").Value = 0.7593360995850622
$ws.Range("S3").Value = 0.1784232365145228
$ws.Range("J4").Value = 0.1449275362318841
$ws.Range("P4").Value = 0.6231884057971014
$ws.Range("S4").Value = 0.2318840579710145
$ws.Range("B6").Value = 0.0705521472392638
$ws.Range("D6").Value = 0.02147239263803681
$ws.Range("F6").Value = 0.0736196319018405
$ws.Range("J6").Value = 0.2822085889570552
$ws.Range("O6").Value = 0.03680981595092025
$ws.Range("Q6").Value = 0.1656441717791411
$ws.Range("R6").Value = 0.05214723926380368
$ws.Range("S6").Value = 0.2975460122699387
$ws.Range("B7").Value = 0.1300448430493273
$ws.Range("D7").Value = 0.01345291479820628
$ws.Range("F7").Value = 0.04035874439461883
$ws.Range("J7").Value = 0.1255605381165919
$ws.Range("O7").Value = 0.04484304932735426
$ws.Range("Q7").Value = 0.1614349775784753
$ws.Range("R7").Value = 0.1210762331838565
$ws.Range("S7").Value = 0.3632286995515695
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.02579365079365079
$ws.Range("F8").Value = 0.06547619047619048
$ws.Range("J8").Value = 0.121031746031746
$ws.Range("O8").Value = 0.02380952380952381
$ws.Range("Q8").Value = 0.2182539682539683
$ws.Range("R8").Value = 0.1170634920634921
$ws.Range("S8").Value = 0.3174603174603174
$ws.Range("B9").Value = 0.09236947791164658
$ws.Range("D9").Value = 0.04819277108433735
$ws.Range("F9").Value = 0.07630522088353414
$ws.Range("J9").Value = 0.1164658634538153
$ws.Range("O9").Value = 0.04016064257028112
$ws.Range("Q9").Value = 0.1847389558232932
$ws.Range("R9").Value = 0.1004016064257028
$ws.Range("S9").Value = 0.3413654618473896
$ws.Range("B10").Value = 0.107993647432504
$ws.Range("D10").Value = 0.01958708311275807
$ws.Range("E10").Value = 0.001058761249338274
$ws.Range("F10").Value = 0.06087877183695077
$ws.Range("J10").Value = 0.1217575436739015
$ws.Range("O10").Value = 0.02541026998411858
$ws.Range("Q10").Value = 0.2535733192165167
$ws.Range("R10").Value = 0.08734780307040763
$ws.Range("S10").Value = 0.3223928004235045
$ws.Range("F11").Value = 0.002652519893899204
$ws.Range("G11").Value = 0.1856763925729443
$ws.Range("J11").Value = 0.1273209549071618
$ws.Range("K11").Value = 0.2307692307692308
$ws.Range("L11").Value = 0.4429708222811671
$ws.Range("S11").Value = 0.01061007957559682
$ws.Range("G12").Value = 0.6763005780346821
$ws.Range("J12").Value = 0.2658959537572254
$ws.Range("K12").Value = 0.01734104046242774
$ws.Range("L12").Value = 0.02890173410404624
$ws.Range("S12").Value = 0.0115606936416185
$ws.Range("G13").Value = 0.5652173913043478
$ws.Range("J13").Value = 0.391304347826087
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02465753424657534
$ws.Range("H15").Value = 0.1479452054794521
$ws.Range("I15").Value = 0.0410958904109589
$ws.Range("J15").Value = 0.389041095890411
$ws.Range("K15").Value = 0.07397260273972603
$ws.Range("M15").Value = 0.005479452054794521
$ws.Range("O15").Value = 0.09315068493150686
$ws.Range("S15").Value = 0.2246575342465753
$ws.Range("F16").Value = 0.02693602693602693
$ws.Range("H16").Value = 0.1279461279461279
$ws.Range("I16").Value = 0.08754208754208755
$ws.Range("J16").Value = 0.4343434343434344
$ws.Range("K16").Value = 0.09427609427609428
$ws.Range("M16").Value = 0.04713804713804714
$ws.Range("O16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.09090909090909091
$ws.Range("F17").Value = 0.03225806451612903
$ws.Range("H17").Value = 0.1626928471248247
$ws.Range("I17").Value = 0.1107994389901823
$ws.Range("J17").Value = 0.4698457223001403
$ws.Range("K17").Value = 0.07854137447405329
$ws.Range("M17").Value = 0.01402524544179523
$ws.Range("N17").Value = 0.002805049088359046
$ws.Range("O17").Value = 0.05890603085553997
$ws.Range("S17").Value = 0.07012622720897616
$ws.Range("F18").Value = 0.02768166089965398
$ws.Range("H18").Value = 0.1349480968858132
$ws.Range("I18").Value = 0.1107266435986159
$ws.Range("J18").Value = 0.4532871972318339
$ws.Range("K18").Value = 0.09342560553633218
$ws.Range("M18").Value = 0.02768166089965398
$ws.Range("O18").Value = 0.07958477508650519
$ws.Range("S18").Value = 0.0726643598615917
$ws.Range("F19").Value = 0.02791696492483894
$ws.Range("H19").Value = 0.1875447387258411
$ws.Range("I19").Value = 0.07015032211882606
$ws.Range("J19").Value = 0.4223335719398711
$ws.Range("K19").Value = 0.09949892627057981
$ws.Range("M19").Value = 0.02720114531138153
$ws.Range("N19").Value = 0.001431639226914817
$ws.Range("O19").Value = 0.08518253400143164
$ws.Range("S19").Value = 0.07874015748031496